$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.728.66'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').Value = '2.925.42'
$ws.Range('E3').Value = '  +0.72%  '
$cell = $ws.Range('D4')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$ws.Range('E4').Value = '  +0.02%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '351.52'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  -0.76%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '106.27'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  -6.81%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -2.29%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '37.69'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  -4.73%  '
$ws.Range('E11').Value = '  +0.96%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '0.0853'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  -2.89%  '
$ws.Range('E13').Value = '  -4.07%  '
$ws.Range('D14').Value = '3.386.69'
$ws.Range('E14').Value = '  +0.67%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '7.66'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -0.91%  '
$ws.Range('D16').Value = '2.934.79'
$ws.Range('E16').Value = '  +1.36%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '0.966'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  -2.08%  '
$ws.Range('D18').Value = '51.649.49'
$ws.Range('E18').Value = '  -0.97%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '3.42'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  +2.22%  '
$ws.Range('E20').Value = '  -3.33%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '13.40'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  -4.79%  '
$ws.Range('D22').Value = '0.0₃0961'
$ws.Range('E22').Value = '  -2.00%  '
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '68.82'
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  -3.13%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '261.94'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  -2.93%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '2.72'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  -2.85%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '0.171'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  -5.05%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '26.50'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  -1.10%  '
$ws.Range('E28').Value = '  +0.19%  '
$ws.Range('E29').Value = '  +8.81%  '
$ws.Range('E30').Value = '  +0.19%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '10.23'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  -3.99%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '35.63'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  -4.63%  '
$ws.Range('E33').Value = '  -5.00%  '
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '5.92'
$cell.Style = 'Normal'
$ws.Range('E34').Value = '  -2.98%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '50.92'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -3.98%  '
$ws.Range('E36').Value = '  -5.67%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('E38').Value = '  -5.64%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '1.96'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  -4.06%  '
$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '17.69'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  -5.89%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '2.66'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  -2.93%  '
$ws.Range('E42').Value = '  -1.00%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '22.60'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  -2.05%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '119.91'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  +1.24%  '
$ws.Range('E45').Value = '  -1.53%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '2.43'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  -3.11%  '
$ws.Range('D47').Value = '2.102.22'
$ws.Range('E47').Value = '  -3.64%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '3.32'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  -6.27%  '
$ws.Range('D49').Value = '3.213.00'
$ws.Range('E49').Value = '  +0.59%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '0.238'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -6.91%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '0.0337'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -5.07%  '
